$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (this shared string is used on Overview!E2/F2/E3/F3 and on the
#    per-language sheets' Status column, C2/C3)
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Column widths widened (report generation bumps these up)
# ---------------------------------------------------------------------------
$overview.Range("E1").ColumnWidth = 29.9777050018311
$overview.Range("F1").ColumnWidth = 29.9777050018311

$zhcn.Range("C1").ColumnWidth = 29.9777050018311
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

$dede.Range("C1").ColumnWidth = 29.9777050018311
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3. Handback data now exists: fill Latest Target File (I), Latest Handback
#    File (J) and Latest Handback DateTime (K) for both data rows on both
#    language sheets. Column I becomes a hyperlink to the source .md file,
#    mirroring column A.
# ---------------------------------------------------------------------------

function Set-HandbackLink($ws, $cellAddr, $targetUrl, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $targetUrl, "", "", $displayText)
    $ws.Range($cellAddr).Font.Underline = 2
    $ws.Range($cellAddr).Font.Color = 15570276
}

$url5b9 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e796d752d3858138bb4613e2ef5dbf02382cc6e/e2e/5b909365-8d13-4ed9-a84f-aa80c853674a.md"
$url6ea = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e796d752d3858138bb4613e2ef5dbf02382cc6e/e2e/6ea5906b-87f6-4695-b6e9-88e0df925214.md"

# zh-cn sheet
Set-HandbackLink $zhcn "I2" $url5b9 "5b909365-8d13-4ed9-a84f-aa80c853674a.md"
$zhcn.Range("J2").Value = "5b909365-8d13-4ed9-a84f-aa80c853674a.311b7811597f92263d7336d1a7899e22e87ef0be.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-19 11:08:24"

Set-HandbackLink $zhcn "I3" $url6ea "6ea5906b-87f6-4695-b6e9-88e0df925214.md"
$zhcn.Range("J3").Value = "6ea5906b-87f6-4695-b6e9-88e0df925214.41eb1c51099641847354ce9af2fd4b40c60ad76f.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-10-19 11:08:24"

# de-de sheet
Set-HandbackLink $dede "I2" $url5b9 "5b909365-8d13-4ed9-a84f-aa80c853674a.md"
$dede.Range("J2").Value = "5b909365-8d13-4ed9-a84f-aa80c853674a.311b7811597f92263d7336d1a7899e22e87ef0be.de-de.xlf"
$dede.Range("K2").Value = "2016-10-19 11:08:41"

Set-HandbackLink $dede "I3" $url6ea "6ea5906b-87f6-4695-b6e9-88e0df925214.md"
$dede.Range("J3").Value = "6ea5906b-87f6-4695-b6e9-88e0df925214.41eb1c51099641847354ce9af2fd4b40c60ad76f.de-de.xlf"
$dede.Range("K3").Value = "2016-10-19 11:08:41"
